# Update "想去人数" (F) / "最低票价" (G) figures across all four sheets
# to match the refreshed scrape (gh-pages output generated at 456a3b4).
$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws = $wb.Worksheets.Item(1)
$ws.Range("F5").Value = 6628
$ws.Range("F6").Value = 513
$ws.Range("F8").Value = 32
$ws.Range("F9").Value = 4544
$ws.Range("F10").Value = 6766
$ws.Range("F11").Value = 1
$ws.Range("F12").Value = 217
$ws.Range("F14").Value = 791
$ws.Range("F15").Value = 109
$ws.Range("F17").Value = 31
$ws.Range("F20").Value = 126
$ws.Range("F22").Value = 184
$ws.Range("F24").Value = 1040
$ws.Range("F25").Value = 506
$ws.Range("F30").Value = 1162
$ws.Range("F32").Value = 85
$ws.Range("G32").Value = 50
$ws.Range("F37").Value = 513
$ws.Range("F38").Value = 347
$ws.Range("F42").Value = 1179
$ws.Range("F43").Value = 517
$ws.Range("F48").Value = 4

# Sheet 2: 演出
$ws = $wb.Worksheets.Item(2)
$ws.Range("F2").Value = 6
$ws.Range("F3").Value = 6
$ws.Range("F12").Value = 120
$ws.Range("F17").Value = 1727
$ws.Range("F22").Value = 189
$ws.Range("F24").Value = 135
$ws.Range("F28").Value = 46
$ws.Range("F31").Value = 735
$ws.Range("F33").Value = 578
$ws.Range("F35").Value = 88
$ws.Range("F37").Value = 9
$ws.Range("F41").Value = 57

# Sheet 3: 本地生活
$ws = $wb.Worksheets.Item(3)
$ws.Range("F4").Value = 709
$ws.Range("F5").Value = 837
$ws.Range("F6").Value = 600
$ws.Range("F8").Value = 1230
$ws.Range("F9").Value = 1102

# Sheet 4: 全部类型
$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 6
$ws.Range("F3").Value = 709
$ws.Range("F6").Value = 837
$ws.Range("F9").Value = 600
$ws.Range("F10").Value = 600
$ws.Range("F12").Value = 6628
$ws.Range("F13").Value = 513
$ws.Range("F15").Value = 32
$ws.Range("F16").Value = 4544
$ws.Range("F18").Value = 6766
$ws.Range("F19").Value = 217
$ws.Range("F22").Value = 791
$ws.Range("F23").Value = 109
$ws.Range("F24").Value = 1230
$ws.Range("F25").Value = 189
$ws.Range("F27").Value = 126
$ws.Range("F28").Value = 184
$ws.Range("F29").Value = 1040
$ws.Range("F31").Value = 509
$ws.Range("F35").Value = 1162
$ws.Range("F36").Value = 85
$ws.Range("G36").Value = 50
$ws.Range("F40").Value = 513
$ws.Range("F41").Value = 578
$ws.Range("F42").Value = 347
$ws.Range("F44").Value = 88
$ws.Range("F46").Value = 517
$ws.Range("F50").Value = 57
